$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at Q (shifts the old "displayString" column Q -> R,
# and the old "tooltip" column R -> S), then populate the new column's
# header and per-row multiplier values.
$ws.Columns("Q:Q").Insert()

$ws.Range("Q1").Value = "valueMultiplier"

$ws.Range("Q2").Value = 1.01
$ws.Range("Q3").Value = 1.01
$ws.Range("Q4").Value = 1.1
$ws.Range("Q5").Value = 1.01
$ws.Range("Q6").Value = 1.1
$ws.Range("Q7").Value = 1.01
$ws.Range("Q8").Value = 1.01
$ws.Range("Q9").Value = 1.1
$ws.Range("Q10").Value = 1.01
$ws.Range("Q11").Value = 1.1
$ws.Range("Q12").Value = 1.25
$ws.Range("Q13").Value = 1.25
$ws.Range("Q14").Value = 1.3
$ws.Range("Q15").Value = 1.05
$ws.Range("Q16").Value = 1.15
$ws.Range("Q17").Value = 1.05
$ws.Range("Q18").Value = 1.15
$ws.Range("Q19").Value = 1.05
$ws.Range("Q20").Value = 1.15
$ws.Range("Q21").Value = 1.05
$ws.Range("Q22").Value = 1.15
$ws.Range("Q23").Value = 1.01
$ws.Range("Q24").Value = 1.1
$ws.Range("Q25").Value = 1.02
$ws.Range("Q26").Value = 1.11
$ws.Range("Q27").Value = 1.2
$ws.Range("Q28").Value = 1.2
$ws.Range("Q29").Value = 1.3
$ws.Range("Q30").Value = 1.02
$ws.Range("Q31").Value = 1.11
$ws.Range("Q32").Value = 1.02
$ws.Range("Q33").Value = 1.11
$ws.Range("Q34").Value = 1.05
$ws.Range("Q35").Value = 1.2

# Match the new column widths (the new valueMultiplier column, and the
# tooltip column that shifted from R to S and got a tighter bestFit width).
$ws.Columns("Q:Q").ColumnWidth = 14.022135416666666
$ws.Columns("S:S").ColumnWidth = 6.166666666666667

# Update the active selection to reflect where the editor ended up.
$ws.Range("Q36").Select()
